# Updates cryptos list values (Price / Volume(1h) columns) to match the
# latest scrape. D-column entries that are valid numeric literals are
# entered with a leading apostrophe so Excel keeps them as text (matching
# the inline-string cell type already used throughout the sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.251.46'
$ws.Range('E2').Value = '  +0.63%  '
$ws.Range('D3').Value = '3.613.46'
$ws.Range('E3').Value = '  +2.25%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('D5').Value = '''601.79'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('D6').Value = '''195.94'
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('E7').Value = '  -0.50%  '
$ws.Range('D8').Value = '''0.999'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +3.15%  '
$ws.Range('D10').Value = '''0.648'
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('D11').Value = '''53.25'
$ws.Range('E11').Value = '  -0.89%  '
$ws.Range('E12').Value = '  +0.54%  '
$ws.Range('D13').Value = '''9.59'
$ws.Range('E13').Value = '  +0.70%  '
$ws.Range('D14').Value = '4.188.57'
$ws.Range('E14').Value = '  +2.18%  '
$ws.Range('D15').Value = '''602.12'
$ws.Range('E15').Value = '  -0.62%  '
$ws.Range('D16').Value = '''12.96'
$ws.Range('E16').Value = '  +2.11%  '
$ws.Range('D17').Value = '70.388.66'
$ws.Range('E17').Value = '  +0.51%  '
$ws.Range('D18').Value = '3.615.49'
$ws.Range('E18').Value = '  +2.01%  '
$ws.Range('D19').Value = '''19.06'
$ws.Range('E19').Value = '  -0.42%  '
$ws.Range('E20').Value = '  +1.65%  '
$ws.Range('E21').Value = '  +0.93%  '
$ws.Range('D22').Value = '''18.59'
$ws.Range('E22').Value = '  +2.38%  '
$ws.Range('D23').Value = '''5.19'
$ws.Range('E23').Value = '  -1.07%  '
$ws.Range('D24').Value = '''103.27'
$ws.Range('E24').Value = '  +0.68%  '
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('D26').Value = '''3.01'
$ws.Range('E26').Value = '  -4.51%  '
$ws.Range('D27').Value = '''10.62'
$ws.Range('E27').Value = '  -2.67%  '
$ws.Range('D28').Value = '''9.70'
$ws.Range('E28').Value = '  +0.79%  '
$ws.Range('D29').Value = '''33.80'
$ws.Range('E29').Value = '  +0.84%  '
$ws.Range('D30').Value = '''4.68'
$ws.Range('E30').Value = '  +7.61%  '
$ws.Range('D31').Value = '''7.28'
$ws.Range('E31').Value = '  +2.48%  '
$ws.Range('E32').Value = '  -2.24%  '
$ws.Range('E33').Value = '  +2.26%  '
$ws.Range('D34').Value = '''63.31'
$ws.Range('E34').Value = '  +0.17%  '
$ws.Range('D35').Value = '0.0₃0883'
$ws.Range('E35').Value = '  +2.94%  '
$ws.Range('D36').Value = '3.934.90'
$ws.Range('E36').Value = '  +5.15%  '
$ws.Range('D37').Value = '''532.67'
$ws.Range('E37').Value = '  +8.99%  '
$ws.Range('E38').Value = '  +0.05%  '
$ws.Range('D39').Value = '''3.05'
$ws.Range('E39').Value = '  +0.22%  '
$ws.Range('D40').Value = '''36.89'
$ws.Range('E40').Value = '  +0.92%  '
$ws.Range('E41').Value = '  -0.86%  '
$ws.Range('E42').Value = '  -2.83%  '
$ws.Range('E43').Value = '  +0.63%  '
$ws.Range('D44').Value = '''0.0461'
$ws.Range('E44').Value = '  +0.55%  '
$ws.Range('D45').Value = '''3.55'
$ws.Range('E45').Value = '  +7.71%  '
$ws.Range('E46').Value = '  +1.43%  '
$ws.Range('D47').Value = '''0.141'
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('E48').Value = '  -0.09%  '
$ws.Range('E49').Value = '  -0.25%  '
$ws.Range('E50').Value = '  -1.63%  '
$ws.Range('E51').Value = '  +1.31%  '
